$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("J6").Value = "FORWARD"
$ws2.Range("J7").Value = "FORWARD"
$ws2.Range("J8").Value = "FORWARD"
$ws2.Range("J9").Value = "FORWARD"
$ws2.Range("J10").Value = "REVERSE"
$ws2.Range("J11").Value = "REVERSE"
$ws2.Range("J12").Value = "REVERSE"
$ws2.Range("J13").Value = "REVERSE"
$ws2.Range("J14").Value = "STOP"
$ws2.Range("J15").Value = "STOP"
$ws2.Range("J16").Value = "STOP"
$ws2.Range("J17").Value = "STOP"
$ws2.Range("J18").Value = "FORWARD"
$ws2.Range("J19").Value = "REVERSE"
$ws2.Range("J20").Value = "REVERSE"
$ws2.Range("J21").Value = "FORWARD"
$ws2.Range("J22").Value = "REVERSE"
$ws2.Range("J23").Value = "FORWARD"
$ws2.Range("J24").Value = "FORWARD"
$ws2.Range("J25").Value = "REVERSE"
$ws2.Range("J26").Value = "REVERSE"
$ws2.Range("J27").Value = "REVERSE"
$ws2.Range("J28").Value = "FORWARD"
$ws2.Range("J29").Value = "FORWARD"
$ws2.Range("J30").Value = "FORWARD"
$ws2.Range("J31").Value = "FORWARD"
$ws2.Range("J32").Value = "REVERSE"
$ws2.Range("J33").Value = "REVERSE"
$ws2.Range("J34").Value = "FORWARD"
$ws2.Range("J35").Value = "STOP"
$ws2.Range("J36").Value = "STOP"
$ws2.Range("J37").Value = "FORWARD"
$ws2.Range("J38").Value = "STOP"
$ws2.Range("J39").Value = "REVERSE"
$ws2.Range("J40").Value = "REVERSE"
$ws2.Range("J41").Value = "STOP"
$ws2.Range("J42").Value = "STOP"
$ws2.Range("J43").Value = "FORWARD"
$ws2.Range("J44").Value = "FORWARD"
$ws2.Range("J45").Value = "STOP"
$ws2.Range("J46").Value = "REVERSE"
$ws2.Range("J47").Value = "STOP"
$ws2.Range("J48").Value = "STOP"
$ws2.Range("J49").Value = "REVERSE"
$ws2.Range("I6").Value = "FRONT_LEFT"
$ws2.Range("I7").Value = "FRONT_RIGHT"
$ws2.Range("I8").Value = "REAR_LEFT"
$ws2.Range("I9").Value = "REAR_RIGHT"
$ws2.Range("I10").Value = "FRONT_LEFT"
$ws2.Range("I11").Value = "FRONT_RIGHT"
$ws2.Range("I12").Value = "REAR_LEFT"
$ws2.Range("I13").Value = "REAR_RIGHT"
$ws2.Range("I14").Value = "FRONT_LEFT"
$ws2.Range("I15").Value = "FRONT_RIGHT"
$ws2.Range("I16").Value = "REAR_LEFT"
$ws2.Range("I17").Value = "REAR_RIGHT"
$ws2.Range("I18").Value = "FRONT_LEFT"
$ws2.Range("I19").Value = "FRONT_RIGHT"
$ws2.Range("I20").Value = "REAR_LEFT"
$ws2.Range("I21").Value = "REAR_RIGHT"
$ws2.Range("I22").Value = "FRONT_LEFT"
$ws2.Range("I23").Value = "FRONT_RIGHT"
$ws2.Range("I24").Value = "REAR_LEFT"
$ws2.Range("I25").Value = "REAR_RIGHT"
$ws2.Range("I26").Value = "FRONT_LEFT"
$ws2.Range("I27").Value = "FRONT_RIGHT"
$ws2.Range("I28").Value = "REAR_LEFT"
$ws2.Range("I29").Value = "REAR_RIGHT"
$ws2.Range("I30").Value = "FRONT_LEFT"
$ws2.Range("I31").Value = "FRONT_RIGHT"
$ws2.Range("I32").Value = "REAR_LEFT"
$ws2.Range("I33").Value = "REAR_RIGHT"
$ws2.Range("I34").Value = "FRONT_LEFT"
$ws2.Range("I35").Value = "FRONT_RIGHT"
$ws2.Range("I36").Value = "REAR_LEFT"
$ws2.Range("I37").Value = "REAR_RIGHT"
$ws2.Range("I38").Value = "FRONT_LEFT"
$ws2.Range("I39").Value = "FRONT_RIGHT"
$ws2.Range("I40").Value = "REAR_LEFT"
$ws2.Range("I41").Value = "REAR_RIGHT"
$ws2.Range("I42").Value = "FRONT_LEFT"
$ws2.Range("I43").Value = "FRONT_RIGHT"
$ws2.Range("I44").Value = "REAR_LEFT"
$ws2.Range("I45").Value = "REAR_RIGHT"
$ws2.Range("I46").Value = "FRONT_LEFT"
$ws2.Range("I47").Value = "FRONT_RIGHT"
$ws2.Range("I48").Value = "REAR_LEFT"
$ws2.Range("I49").Value = "REAR_RIGHT"
$ws2.Range("H6").Value = "forward"
$ws2.Range("H10").Value = "backward"
$ws2.Range("H14").Value = "stopped"
$ws2.Range("H18").Value = "strafeRight"
$ws2.Range("H22").Value = "strafeLeft"
$ws2.Range("H26").Value = "rotateRight"
$ws2.Range("H30").Value = "rotateLeft"
$ws2.Range("H34").Value = "rightForward"
$ws2.Range("H38").Value = "rightBackward"
$ws2.Range("H42").Value = "leftForward"
$ws2.Range("H46").Value = "leftBackward"
$ws2.Range("C3").Value = "LEFT SIDE"
$ws2.Range("E3").Value = "RIGHT SIDE"
$ws2.Range("B4").Value = "direction"
$ws2.Range("C4").Value = "Front"
$ws2.Range("D4").Value = "Back"
$ws2.Range("E4").Value = "Front"
$ws2.Range("F4").Value = "Back"

[void]$ws1.Range("A2:E14").Select()
[void]$ws2.Range("I46:I49").Select()
[void]$ws2.Activate()
